$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 4.793498
$ws.Cells.Item(2, 8).Value = 14.380494
$ws.Cells.Item(2, 9).Value = 0.6565493962715211
$ws.Cells.Item(2, 10).Value = 0.6565493962715211
$ws.Cells.Item(2, 13).Value = 0.2054816666666667
$ws.Cells.Item(2, 14).Value = 0.616445
$ws.Cells.Item(2, 15).Value = 0.00496424614546655
$ws.Cells.Item(2, 16).Value = 0.004964246145466549
$ws.Cells.Item(2, 17).Value = 0.9849759582033335
$ws.Cells.Item(2, 18).Value = 8.86478362383
$ws.Cells.Item(2, 19).Value = 0.003259272809749289
$ws.Cells.Item(2, 20).Value = 0.003259272809749288

# Row 3
$ws.Cells.Item(3, 7).Value = 4.793498
$ws.Cells.Item(3, 8).Value = 14.380494
$ws.Cells.Item(3, 9).Value = 0.6565493962715211
$ws.Cells.Item(3, 10).Value = 0.6565493962715211
$ws.Cells.Item(3, 15).Value = 0.9529850468799925
$ws.Cells.Item(3, 16).Value = 0.9529850468799924
$ws.Cells.Item(3, 17).Value = 189.0855796023087
$ws.Cells.Item(3, 18).Value = 1701.770216420778
$ws.Cells.Item(3, 19).Value = 0.6256817571848463
$ws.Cells.Item(3, 20).Value = 0.6256817571848462

# Row 4
$ws.Cells.Item(4, 7).Value = 4.793498
$ws.Cells.Item(4, 8).Value = 14.380494
$ws.Cells.Item(4, 9).Value = 0.6565493962715211
$ws.Cells.Item(4, 10).Value = 0.6565493962715211
$ws.Cells.Item(4, 13).Value = 1.712817
$ws.Cells.Item(4, 14).Value = 5.138451
$ws.Cells.Item(4, 15).Value = 0.04138006727350978
$ws.Cells.Item(4, 16).Value = 0.04138006727350978
$ws.Cells.Item(4, 17).Value = 8.210384863866
$ws.Cells.Item(4, 18).Value = 73.893463774794
$ws.Cells.Item(4, 19).Value = 0.02716805818609778
$ws.Cells.Item(4, 20).Value = 0.02716805818609777

# Row 5
$ws.Cells.Item(5, 7).Value = 4.793498
$ws.Cells.Item(5, 8).Value = 14.380494
$ws.Cells.Item(5, 9).Value = 0.6565493962715211
$ws.Cells.Item(5, 10).Value = 0.6565493962715211
$ws.Cells.Item(5, 13).Value = 0.02775933333333333
$ws.Cells.Item(5, 14).Value = 0.083278
$ws.Cells.Item(5, 15).Value = 0.000670639701031176
$ws.Cells.Item(5, 16).Value = 0.000670639701031176
$ws.Cells.Item(5, 17).Value = 0.1330643088146667
$ws.Cells.Item(5, 18).Value = 1.197578779332
$ws.Cells.Item(5, 19).Value = 0.000440308090827732
$ws.Cells.Item(5, 20).Value = 0.000440308090827732

# Row 6
$ws.Cells.Item(6, 9).Value = 0.02833032029515766
$ws.Cells.Item(6, 10).Value = 0.02833032029515767
$ws.Cells.Item(6, 13).Value = 0.2054816666666667
$ws.Cells.Item(6, 14).Value = 0.616445
$ws.Cells.Item(6, 15).Value = 0.00496424614546655
$ws.Cells.Item(6, 16).Value = 0.004964246145466549
$ws.Cells.Item(6, 17).Value = 0.04250203341499999
$ws.Cells.Item(6, 18).Value = 0.382518300735
$ws.Cells.Item(6, 19).Value = 0.0001406386833250692
$ws.Cells.Item(6, 20).Value = 0.0001406386833250692

# Row 7
$ws.Cells.Item(7, 9).Value = 0.02833032029515766
$ws.Cells.Item(7, 10).Value = 0.02833032029515767
$ws.Cells.Item(7, 15).Value = 0.9529850468799925
$ws.Cells.Item(7, 16).Value = 0.9529850468799924
$ws.Cells.Item(7, 19).Value = 0.02699837161460603
$ws.Cells.Item(7, 20).Value = 0.02699837161460603

# Row 8
$ws.Cells.Item(8, 9).Value = 0.02833032029515766
$ws.Cells.Item(8, 10).Value = 0.02833032029515767
$ws.Cells.Item(8, 13).Value = 1.712817
$ws.Cells.Item(8, 14).Value = 5.138451
$ws.Cells.Item(8, 15).Value = 0.04138006727350978
$ws.Cells.Item(8, 16).Value = 0.04138006727350978
$ws.Cells.Item(8, 17).Value = 0.354280781097
$ws.Cells.Item(8, 18).Value = 3.188527029873
$ws.Cells.Item(8, 19).Value = 0.001172310559693704
$ws.Cells.Item(8, 20).Value = 0.001172310559693704

# Row 9
$ws.Cells.Item(9, 9).Value = 0.02833032029515766
$ws.Cells.Item(9, 10).Value = 0.02833032029515767
$ws.Cells.Item(9, 13).Value = 0.02775933333333333
$ws.Cells.Item(9, 14).Value = 0.083278
$ws.Cells.Item(9, 15).Value = 0.000670639701031176
$ws.Cells.Item(9, 16).Value = 0.000670639701031176
$ws.Cells.Item(9, 17).Value = 0.005741768265999999
$ws.Cells.Item(9, 18).Value = 0.051675914394
$ws.Cells.Item(9, 19).Value = 0.00001899943753286199
$ws.Cells.Item(9, 20).Value = 0.000018999437532862

# Row 10
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.147564
$ws.Cells.Item(10, 8).Value = 0.442692
$ws.Cells.Item(10, 9).Value = 0.02021134776971029
$ws.Cells.Item(10, 10).Value = 0.02021134776971029
$ws.Cells.Item(10, 13).Value = 0.2054816666666667
$ws.Cells.Item(10, 14).Value = 0.616445
$ws.Cells.Item(10, 15).Value = 0.00496424614546655
$ws.Cells.Item(10, 16).Value = 0.004964246145466549
$ws.Cells.Item(10, 17).Value = 0.03032169666
$ws.Cells.Item(10, 18).Value = 0.27289526994
$ws.Cells.Item(10, 19).Value = 0.0001003341052604682
$ws.Cells.Item(10, 20).Value = 0.0001003341052604682

# Row 11
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.147564
$ws.Cells.Item(11, 8).Value = 0.442692
$ws.Cells.Item(11, 9).Value = 0.02021134776971029
$ws.Cells.Item(11, 10).Value = 0.02021134776971029
$ws.Cells.Item(11, 15).Value = 0.9529850468799925
$ws.Cells.Item(11, 16).Value = 0.9529850468799924
$ws.Cells.Item(11, 17).Value = 5.820848254955999
$ws.Cells.Item(11, 18).Value = 52.387634294604
$ws.Cells.Item(11, 19).Value = 0.01926111220182519
$ws.Cells.Item(11, 20).Value = 0.01926111220182519

# Row 12
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.147564
$ws.Cells.Item(12, 8).Value = 0.442692
$ws.Cells.Item(12, 9).Value = 0.02021134776971029
$ws.Cells.Item(12, 10).Value = 0.02021134776971029
$ws.Cells.Item(12, 13).Value = 1.712817
$ws.Cells.Item(12, 14).Value = 5.138451
$ws.Cells.Item(12, 15).Value = 0.04138006727350978
$ws.Cells.Item(12, 16).Value = 0.04138006727350978
$ws.Cells.Item(12, 17).Value = 0.252750127788
$ws.Cells.Item(12, 18).Value = 2.274751150092
$ws.Cells.Item(12, 19).Value = 0.0008363469303989137
$ws.Cells.Item(12, 20).Value = 0.0008363469303989136

# Row 13
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.147564
$ws.Cells.Item(13, 8).Value = 0.442692
$ws.Cells.Item(13, 9).Value = 0.02021134776971029
$ws.Cells.Item(13, 10).Value = 0.02021134776971029
$ws.Cells.Item(13, 13).Value = 0.02775933333333333
$ws.Cells.Item(13, 14).Value = 0.083278
$ws.Cells.Item(13, 15).Value = 0.000670639701031176
$ws.Cells.Item(13, 16).Value = 0.000670639701031176
$ws.Cells.Item(13, 17).Value = 0.004096278264
$ws.Cells.Item(13, 18).Value = 0.036866504376
$ws.Cells.Item(13, 19).Value = 0.00001355453222571564
$ws.Cells.Item(13, 20).Value = 0.00001355453222571564

# Row 14
$ws.Cells.Item(14, 7).Value = 2.153144
$ws.Cells.Item(14, 8).Value = 6.459432
$ws.Cells.Item(14, 9).Value = 0.2949089356636109
$ws.Cells.Item(14, 10).Value = 0.294908935663611
$ws.Cells.Item(14, 13).Value = 0.2054816666666667
$ws.Cells.Item(14, 14).Value = 0.616445
$ws.Cells.Item(14, 15).Value = 0.00496424614546655
$ws.Cells.Item(14, 16).Value = 0.004964246145466549
$ws.Cells.Item(14, 17).Value = 0.4424316176933333
$ws.Cells.Item(14, 18).Value = 3.98188455924
$ws.Cells.Item(14, 19).Value = 0.001464000547131723
$ws.Cells.Item(14, 20).Value = 0.001464000547131723

# Row 15
$ws.Cells.Item(15, 7).Value = 2.153144
$ws.Cells.Item(15, 8).Value = 6.459432
$ws.Cells.Item(15, 9).Value = 0.2949089356636109
$ws.Cells.Item(15, 10).Value = 0.294908935663611
$ws.Cells.Item(15, 15).Value = 0.9529850468799925
$ws.Cells.Item(15, 16).Value = 0.9529850468799924
$ws.Cells.Item(15, 17).Value = 84.93348306544264
$ws.Cells.Item(15, 18).Value = 764.401347588984
$ws.Cells.Item(15, 19).Value = 0.281043805878715
$ws.Cells.Item(15, 20).Value = 0.281043805878715

# Row 16
$ws.Cells.Item(16, 7).Value = 2.153144
$ws.Cells.Item(16, 8).Value = 6.459432
$ws.Cells.Item(16, 9).Value = 0.2949089356636109
$ws.Cells.Item(16, 10).Value = 0.294908935663611
$ws.Cells.Item(16, 13).Value = 1.712817
$ws.Cells.Item(16, 14).Value = 5.138451
$ws.Cells.Item(16, 15).Value = 0.04138006727350978
$ws.Cells.Item(16, 16).Value = 0.04138006727350978
$ws.Cells.Item(16, 17).Value = 3.687941646648
$ws.Cells.Item(16, 18).Value = 33.191474819832
$ws.Cells.Item(16, 19).Value = 0.01220335159731939
$ws.Cells.Item(16, 20).Value = 0.01220335159731939

# Row 17
$ws.Cells.Item(17, 7).Value = 2.153144
$ws.Cells.Item(17, 8).Value = 6.459432
$ws.Cells.Item(17, 9).Value = 0.2949089356636109
$ws.Cells.Item(17, 10).Value = 0.294908935663611
$ws.Cells.Item(17, 13).Value = 0.02775933333333333
$ws.Cells.Item(17, 14).Value = 0.083278
$ws.Cells.Item(17, 15).Value = 0.000670639701031176
$ws.Cells.Item(17, 16).Value = 0.000670639701031176
$ws.Cells.Item(17, 17).Value = 0.05976984201066666
$ws.Cells.Item(17, 18).Value = 0.537928578096
$ws.Cells.Item(17, 19).Value = 0.0001977776404448664
$ws.Cells.Item(17, 20).Value = 0.0001977776404448664
